# Update the "Synonym" worksheet with 5 new gene/synonym rows, inserted in their
# correct alphabetically-sorted position (the sheet is kept sorted by column C).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Synonym")

# Each entry: (final row number after all inserts, Gene synonym in C, actual gene name in D)
# Inserting in increasing order of final row number works because every insertion
# position is specified in terms of the final (already-shifted) layout.
$ws.Rows.Item(7).Insert()
$ws.Range("C7").Value = "ACKR1"
$ws.Range("D7").Value = "DARC"

$ws.Rows.Item(41).Insert()
$ws.Range("C41").Value = "JCHAIN"
$ws.Range("D41").Value = "IGJ"

$ws.Rows.Item(60).Insert()
$ws.Range("C60").Value = "TBC1D31"
$ws.Range("D60").Value = "WDR67"

$ws.Rows.Item(64).Insert()
$ws.Range("C64").Value = "VEGFD"
$ws.Range("D64").Value = "FIGF"

$ws.Rows.Item(67).Insert()
$ws.Range("C67").Value = "XCL2"
$ws.Range("D67").Value = "SCYC2"

# Restore the sheet's selection/active-cell and scroll position to match the
# post-edit state (user ended up having just added the last row and selected C65).
$ws.Activate()
$ws.Application.Goto($ws.Range("A49"))
$ws.Range("C65").Select()
